$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.820.04'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.476.48'
$ws.Range("E3").Value = '  -2.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.72'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.63'
$ws.Range("E6").Value = '  -4.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.552'
$ws.Range("E7").Value = '  -2.55%  '
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.507'
$ws.Range("E9").Value = '  -3.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.94'
$ws.Range("E10").Value = '  -4.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0783'
$ws.Range("E11").Value = '  -2.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.108'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.98'
$ws.Range("E13").Value = '  -4.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.871.37'
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.439.87'
$ws.Range("E15").Value = '  -3.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.71'
$ws.Range("E16").Value = '  -6.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.788'
$ws.Range("E17").Value = '  -4.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.845.42'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").Value = '  -5.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0918'
$ws.Range("E20").Value = '  -3.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.65'
$ws.Range("E21").Value = '  -4.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.34'
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.28'
$ws.Range("E23").Value = '  -3.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.79'
$ws.Range("E24").Value = '  -3.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.93'
$ws.Range("E25").Value = '  -5.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.72'
$ws.Range("E27").Value = '  -4.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").Value = '  -4.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.73'
$ws.Range("E29").Value = '  -3.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.46'
$ws.Range("E30").Value = '  -7.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '154.59'
$ws.Range("E31").Value = '  -1.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.60'
$ws.Range("E32").Value = '  -2.07%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0756'
$ws.Range("E34").Value = '  -4.91%  '
$ws.Range("B35").Value = 'ApeXProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.56'
$ws.Range("E35").Value = '  -8.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.03'
$ws.Range("E36").Value = '  -4.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.14'
$ws.Range("E37").Value = '  -5.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.88'
$ws.Range("E38").Value = '  -6.87%  '
$ws.Range("E39").Value = '  -3.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.114'
$ws.Range("E40").Value = '  -3.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.01'
$ws.Range("E41").Value = '  -6.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.99'
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.996.33'
$ws.Range("E44").Value = '  +1.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0286'
$ws.Range("E45").Value = '  -3.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.07'
$ws.Range("E46").Value = '  -7.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.67'
$ws.Range("E47").Value = '  -2.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.725.79'
$ws.Range("E48").Value = '  -1.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '76.52'
$ws.Range("E49").Value = '  -5.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.181'
$ws.Range("E50").Value = '  -5.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.61'
$ws.Range("E51").Value = '  -3.96%  '
